$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.409.51'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '3.689.35'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'679.85"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").Value = "'160.80"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'0.495"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").Value = "'7.17"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").Value = "'0.439"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +0.48%  '

$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '4.310.36'
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = "'32.48"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -0.12%  '

$ws.Range("D15").Value = '3.699.57'
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("D16").Value = '69.373.93'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = "'0.117"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +2.71%  '

$ws.Range("D18").Value = "'16.04"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("D19").Value = "'6.47"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").Value = "'472.10"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.59%  '

$ws.Range("D21").Value = "'9.83"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -0.72%  '

$ws.Range("D22").Value = "'0.650"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").Value = "'80.27"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("D24").Value = '3.836.40'
$ws.Range("E24").Value = '  +0.14%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").Value = "'10.91"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -1.09%  '

$ws.Range("D28").Value = "'9.14"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.63%  '

$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("D30").Value = "'1.74"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("D31").Value = "'2.02"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.46%  '

$ws.Range("D32").Value = "'6.58"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -1.30%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("D34").Value = "'26.99"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("D35").Value = '3.678.58'
$ws.Range("E35").Value = '  +0.49%  '

$ws.Range("E36").Value = '  +1.81%  '

$ws.Range("D37").Value = "'8.48"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +3.34%  '

$ws.Range("D38").Value = "'6.22"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D40").Value = "'2.26"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").Value = "'0.0904"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").Value = "'168.58"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +1.90%  '

$ws.Range("D44").Value = "'0.941"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").Value = "'46.71"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.33%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = "'2.73"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("B47").Value = 'FLOKI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D47").Value = "'0.000280"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = "'28.01"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.22%  '

$ws.Range("D49").Value = "'1.29"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -1.90%  '

$ws.Range("E50").Value = '  -2.22%  '

$ws.Range("D51").Value = "'7.88"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +0.30%  '
